$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous used range completely so stale column C/D values are removed
$ws.Range("A1:D68").ClearContents()

$ws.Cells.Item(1,1).Value = "<Darkness:250>"
$ws.Cells.Item(1,2).Value = "<Darkness:250>"
$ws.Cells.Item(2,1).Value = "自動実行"
$ws.Cells.Item(2,2).Value = "自動実行"
$ws.Cells.Item(3,1).Value = "SetCharLight -1 2 1 80 0.2"
$ws.Cells.Item(3,2).Value = "SetCharLight -1 2 1 80 0.2"
$ws.Cells.Item(4,1).Value = "燭台"
$ws.Cells.Item(4,2).Value = "燭台"
$ws.Cells.Item(5,1).Value = " <enemy:99>"
$ws.Cells.Item(5,2).Value = " <enemy:99>"
$ws.Cells.Item(6,1).Value = "燭台だ。`n固定されて動かせないが灯りに使えるかもしれない。"
$ws.Cells.Item(6,2).Value = "燭台だ。`n固定されて動かせないが灯りに使えるかもしれない。"
$ws.Cells.Item(7,1).Value = "変数203（ARGP攻撃種類）`n1斬　2打撃　3水　4火　5雷`n6誘惑　7食べ物`n特殊206（個別攻撃種類）`n1リンゴ　2皿　3卵"
$ws.Cells.Item(8,1).Value = "灯り 6 1 200 0.1"
$ws.Cells.Item(9,1).Value = "捕まえリリー"
$ws.Cells.Item(9,2).Value = "捕まえリリー"
$ws.Cells.Item(10,1).Value = "\n<シィナ>カワバンガにゃ！`n待ち伏せ成功にゃ！"
$ws.Cells.Item(10,2).Value = "\n<シィナ>カワバンガにゃ！`n待ち伏せ成功にゃ！"
$ws.Cells.Item(11,1).Value = "\n<シィナ>もう逃がさんにゃぁ♥`nじゅるり♥"
$ws.Cells.Item(11,2).Value = "\n<シィナ>もう逃がさんにゃぁ♥`nじゅるり♥"
$ws.Cells.Item(12,1).Value = "\n<シィナ>アタシに搾って欲しいのにゃ？`nしょうがないにゃぁ・・・"
$ws.Cells.Item(12,2).Value = "\n<シィナ>アタシに搾って欲しいのにゃ？`nしょうがないにゃぁ・・・"
$ws.Cells.Item(13,1).Value = " <enemy:99><cw:1.5>"
$ws.Cells.Item(13,2).Value = " <enemy:99><cw:1.5>"
$ws.Cells.Item(14,1).Value = "ビン"
$ws.Cells.Item(14,2).Value = "ビン"
$ws.Cells.Item(15,1).Value = "割れたビンを踏んでしまった・・・！"
$ws.Cells.Item(15,2).Value = "割れたビンを踏んでしまった・・・！"
$ws.Cells.Item(16,1).Value = "\n<シィナ>ん？`nにゃるほど。そこにいたにゃ。"
$ws.Cells.Item(16,2).Value = "\n<シィナ>ん？`nにゃるほど。そこにいたにゃ。"
$ws.Cells.Item(17,1).Value = "MP_SET_MOVIE e14_Ta"
$ws.Cells.Item(17,2).Value = "MP_SET_MOVIE e14_Ta"
$ws.Cells.Item(18,1).Value = "MP_SET_LOOP 80 on"
$ws.Cells.Item(18,2).Value = "MP_SET_LOOP 80 on"
$ws.Cells.Item(19,1).Value = "\n<シィナ>袋小路だし♥`nもう絶対逃げられんにゃ。`nほーら。食ってやるにゃん♥"
$ws.Cells.Item(19,2).Value = "\n<シィナ>袋小路だし♥`nもう絶対逃げられんにゃ。`nほーら。食ってやるにゃん♥"
$ws.Cells.Item(20,1).Value = "ビン踏みシィナ"
$ws.Cells.Item(20,2).Value = "ビン踏みシィナ"
$ws.Cells.Item(21,1).Value = "灯り 2 1 120 0.1"
$ws.Cells.Item(22,1).Value = "\n<シィナ>追い込んだにゃ♥`nほーら捕まえた♥"
$ws.Cells.Item(22,2).Value = "\n<シィナ>追い込んだにゃ♥`nほーら捕まえた♥"
$ws.Cells.Item(23,1).Value = "EV021"
$ws.Cells.Item(23,2).Value = "EV021"
$ws.Cells.Item(24,1).Value = "水たまりを踏んでしまった・・・"
$ws.Cells.Item(24,2).Value = "水たまりを踏んでしまった・・・"
$ws.Cells.Item(25,1).Value = "\n<ライム>私だよー！`n待ち伏せ成功ー♪"
$ws.Cells.Item(25,2).Value = "\n<ライム>私だよー！`n待ち伏せ成功ー♪"
$ws.Cells.Item(26,1).Value = "EV022"
$ws.Cells.Item(26,2).Value = "EV022"
$ws.Cells.Item(27,1).Value = "\n<リリー>うわぁ！！びっくりしたぁ！！"
$ws.Cells.Item(27,2).Value = "\n<リリー>うわぁ！！びっくりしたぁ！！"
$ws.Cells.Item(28,1).Value = "\n<リリー>あ、あんたか・・・`nま、待ち伏せ成功ね！`nおほほほほ♥（心臓止まるかと思った・・・）"
$ws.Cells.Item(28,2).Value = "\n<リリー>あ、あんたか・・・`nま、待ち伏せ成功ね！`nおほほほほ♥（心臓止まるかと思った・・・）"
$ws.Cells.Item(29,1).Value = "\n<リリー>さぁ、捕まえたわよ♥`nもう逃がさないから。"
$ws.Cells.Item(29,2).Value = "\n<リリー>さぁ、捕まえたわよ♥`nもう逃がさないから。"
$ws.Cells.Item(30,1).Value = "\n<リリー>見えてるのに。`n自分から搾られに来たの？`nくす♥"
$ws.Cells.Item(30,2).Value = "\n<リリー>見えてるのに。`n自分から搾られに来たの？`nくす♥"
$ws.Cells.Item(31,1).Value = "捕まえライム"
$ws.Cells.Item(31,2).Value = "捕まえライム"
$ws.Cells.Item(32,1).Value = "\n<ライム>待ち伏せ成功ー♥`nにゅるりぃー♥"
$ws.Cells.Item(32,2).Value = "\n<ライム>待ち伏せ成功ー♥`nにゅるりぃー♥"
$ws.Cells.Item(33,1).Value = "\n<ライム>真っ暗の中で待ってた甲斐があったよー♥`nもう逃げられないからねー♥`nえーい♥"
$ws.Cells.Item(33,2).Value = "\n<ライム>真っ暗の中で待ってた甲斐があったよー♥`nもう逃げられないからねー♥`nえーい♥"
$ws.Cells.Item(34,1).Value = "\n<ライム>あっ♥自分から捕まりに来た♥`nえへへー♥`nそっかそっかー♥"
$ws.Cells.Item(34,2).Value = "\n<ライム>あっ♥自分から捕まりに来た♥`nえへへー♥`nそっかそっかー♥"
$ws.Cells.Item(35,1).Value = "\n<ライム>パイズリして欲しくなったんでしょー。`n分かったよー♥`nほら、おちんちんちょーだい♥"
$ws.Cells.Item(35,2).Value = "\n<ライム>パイズリして欲しくなったんでしょー。`n分かったよー♥`nほら、おちんちんちょーだい♥"
$ws.Cells.Item(36,1).Value = "EV025"
$ws.Cells.Item(36,2).Value = "EV025"
$ws.Cells.Item(37,1).Value = "食糧庫ドア"
$ws.Cells.Item(37,2).Value = "食糧庫ドア"
$ws.Cells.Item(38,1).Value = "\n<シィナ>・・・`nあいつ全然来んにゃぁ。"
$ws.Cells.Item(38,2).Value = "\n<シィナ>・・・`nあいつ全然来んにゃぁ。"
$ws.Cells.Item(39,1).Value = "\n<ライム>倉庫に隠れてたらロープ目当てに来ると`n思ってたんだけどなー。"
$ws.Cells.Item(39,2).Value = "\n<ライム>倉庫に隠れてたらロープ目当てに来ると`n思ってたんだけどなー。"
$ws.Cells.Item(40,1).Value = "\n<リリー>ぐぬぬ・・・`nこんな暗闇で息潜めてバカみたいじゃない。`n早く来なさいよ。"
$ws.Cells.Item(40,2).Value = "\n<リリー>ぐぬぬ・・・`nこんな暗闇で息潜めてバカみたいじゃない。`n早く来なさいよ。"
$ws.Cells.Item(41,1).Value = "\n<シィナ>にゃ！`nロープが無くなってるにゃん！"
$ws.Cells.Item(41,2).Value = "\n<シィナ>にゃ！`nロープが無くなってるにゃん！"
$ws.Cells.Item(42,1).Value = "\n<ライム>ウソ！`nこんなに暗いのにどうやってー？"
$ws.Cells.Item(42,2).Value = "\n<ライム>ウソ！`nこんなに暗いのにどうやってー？"
$ws.Cells.Item(43,1).Value = "\n<リリー>す、凄みよ・・・！`nあ・・・あいつ凄みで！`n倉庫のロープを探知したんだわ・・・！"
$ws.Cells.Item(43,2).Value = "\n<リリー>す、凄みよ・・・！`nあ・・・あいつ凄みで！`n倉庫のロープを探知したんだわ・・・！"
$ws.Cells.Item(44,1).Value = "\n<ライム>凄み！？"
$ws.Cells.Item(44,2).Value = "\n<ライム>凄み！？"
$ws.Cells.Item(45,1).Value = "\n<シィナ>にゃぁぁ・・・`n凄みかぁ・・・"
$ws.Cells.Item(45,2).Value = "\n<シィナ>にゃぁぁ・・・`n凄みかぁ・・・"
$ws.Cells.Item(46,1).Value = "\n[1]は\C[3]『凄み』\C[0]を覚えた！"
$ws.Cells.Item(46,2).Value = "\n[1]は\C[3]『凄み』\C[0]を覚えた！"
$ws.Cells.Item(47,1).Value = "ロープ"
$ws.Cells.Item(47,2).Value = "Rope"
$ws.Cells.Item(48,1).Value = "ロープ・・・！`nこれは何かの役に立ちそうだ・・・！"
$ws.Cells.Item(48,2).Value = "ロープ・・・！`nこれは何かの役に立ちそうだ・・・！"
$ws.Cells.Item(49,1).Value = "EV031"
$ws.Cells.Item(49,2).Value = "EV031"
$ws.Cells.Item(50,1).Value = "罠を踏んでしまった！！`n身体が痺れて動かない・・・！"
$ws.Cells.Item(50,2).Value = "罠を踏んでしまった！！`n身体が痺れて動かない・・・！"
$ws.Cells.Item(51,1).Value = "\n<リリー>捕獲成功～♪`nこんな暗い中、足元をよく見ずに歩いてるからよ。`nばーか♥"
$ws.Cells.Item(51,2).Value = "\n<リリー>捕獲成功～♪`nこんな暗い中、足元をよく見ずに歩いてるからよ。`nばーか♥"
$ws.Cells.Item(52,1).Value = "\n<リリー>灯り付けてたのに。`n罠が見えなかったの？`nくす♥"
$ws.Cells.Item(52,2).Value = "\n<リリー>灯り付けてたのに。`n罠が見えなかったの？`nくす♥"
$ws.Cells.Item(53,1).Value = "EV032"
$ws.Cells.Item(53,2).Value = "EV032"
$ws.Cells.Item(54,1).Value = "立体起動"
$ws.Cells.Item(54,2).Value = "立体起動"
$ws.Cells.Item(55,1).Value = "<enemy:99><cw:1.3><ch:1.5>"
$ws.Cells.Item(55,2).Value = "<enemy:99><cw:1.3><ch:1.5>"
$ws.Cells.Item(56,1).Value = "箱のようなものがあるが`n高くて手が届かない・・・"
$ws.Cells.Item(56,2).Value = "箱のようなものがあるが`n高くて手が届かない・・・"
$ws.Cells.Item(57,1).Value = "リリーの日記"
$ws.Cells.Item(57,2).Value = "リリーの日記"
$ws.Cells.Item(58,1).Value = "ライムの日記`n見られているような・・・"
$ws.Cells.Item(58,2).Value = "ライムの日記`n見られているような・・・"
$ws.Cells.Item(59,1).Value = "読んでみる"
$ws.Cells.Item(59,2).Value = "読んでみる"
$ws.Cells.Item(60,1).Value = "やめておく"
$ws.Cells.Item(60,2).Value = "やめておく"
$ws.Cells.Item(61,1).Value = "何かあるが暗くてよく見えない・・・"
$ws.Cells.Item(61,2).Value = "何かあるが暗くてよく見えない・・・"
$ws.Cells.Item(62,1).Value = "リリーの日記`nオークはどこから"
$ws.Cells.Item(62,2).Value = "リリーの日記`nオークはどこから"
$ws.Cells.Item(63,1).Value = "EV037"
$ws.Cells.Item(63,2).Value = "EV037"
$ws.Cells.Item(64,1).Value = "EV038"
$ws.Cells.Item(64,2).Value = "EV038"
$ws.Cells.Item(65,1).Value = "EV039"
$ws.Cells.Item(65,2).Value = "EV039"
$ws.Cells.Item(66,1).Value = "EV040"
$ws.Cells.Item(66,2).Value = "EV040"
$ws.Cells.Item(67,1).Value = "EV041"
$ws.Cells.Item(67,2).Value = "EV041"
$ws.Cells.Item(68,1).Value = "灯りを付けるには・・・"
$ws.Cells.Item(68,2).Value = "灯りを付けるには・・・"
$ws.Cells.Item(69,1).Value = "答えを見る"
$ws.Cells.Item(69,2).Value = "答えを見る"
$ws.Cells.Item(70,1).Value = "見ない"
$ws.Cells.Item(70,2).Value = "見ない"
$ws.Cells.Item(71,1).Value = "炎系の攻撃で燭台に火を着けてください。"
$ws.Cells.Item(71,2).Value = "炎系の攻撃で燭台に火を着けてください。"

# Values containing embedded newlines trigger an auto row-height bump in this
# COM host (mimicking Excel's wrap-driven autofit). The target canonical XML
# has no custom row heights, so reset every touched row back to standard.
$ws.Range("A1:B71").EntireRow.AutoFit()
